$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Row 2
Set-TextValue $ws.Range("D2") '36.566.11'
Set-TextValue $ws.Range("E2") '  -0.17%  '

# Row 3
Set-TextValue $ws.Range("D3") '1.958.83'
Set-TextValue $ws.Range("E3") '  -0.23%  '

# Row 4
Set-TextValue $ws.Range("E4") '  -0.04%  '

# Row 5
Set-TextValue $ws.Range("D5") '244.54'
Set-TextValue $ws.Range("E5") '  -0.34%  '

# Row 6
Set-TextValue $ws.Range("E6") '  -0.72%  '

# Row 7
Set-TextValue $ws.Range("D7") '58.73'
Set-TextValue $ws.Range("E7") '  -1.36%  '

# Row 8
Set-TextValue $ws.Range("E8") '  -0.05%  '

# Row 9
Set-TextValue $ws.Range("E9") '  +2.21%  '

# Row 10
Set-TextValue $ws.Range("D10") '0.0806'
Set-TextValue $ws.Range("E10") '  -6.32%  '

# Row 11
Set-TextValue $ws.Range("E11") '  -0.86%  '

# Row 12
Set-TextValue $ws.Range("D12") '22.11'
Set-TextValue $ws.Range("E12") '  -2.21%  '

# Row 13
Set-TextValue $ws.Range("D13") '0.831'
Set-TextValue $ws.Range("E13") '  -1.04%  '

# Row 14
Set-TextValue $ws.Range("D14") '2.246.19'
Set-TextValue $ws.Range("E14") '  -0.26%  '

# Row 15
Set-TextValue $ws.Range("D15") '13.71'
Set-TextValue $ws.Range("E15") '  -0.96%  '

# Row 16
Set-TextValue $ws.Range("E16") '  -0.22%  '

# Row 17
Set-TextValue $ws.Range("D17") '1.947.21'
Set-TextValue $ws.Range("E17") '  -0.40%  '

# Row 18
Set-TextValue $ws.Range("D18") '36.468.91'
Set-TextValue $ws.Range("E18") '  -0.39%  '

# Row 19
Set-TextValue $ws.Range("E19") '  -0.94%  '

# Row 20
Set-TextValue $ws.Range("D20") '0.0₃0854'
Set-TextValue $ws.Range("E20") '  -2.73%  '

# Row 21
Set-TextValue $ws.Range("D21") '228.55'
Set-TextValue $ws.Range("E21") '  -1.12%  '

# Row 22
Set-TextValue $ws.Range("D22") '5.05'
Set-TextValue $ws.Range("E22") '  -1.25%  '

# Row 23
Set-TextValue $ws.Range("E23") '  -0.08%  '

# Row 24
Set-TextValue $ws.Range("D24") '2.45'
Set-TextValue $ws.Range("E24") '  -1.48%  '

# Row 25
Set-TextValue $ws.Range("E25") '  +1.24%  '

# Row 26
Set-TextValue $ws.Range("D26") '9.25'
Set-TextValue $ws.Range("E26") '  -2.23%  '

# Row 27
Set-TextValue $ws.Range("D27") '0.138'
Set-TextValue $ws.Range("E27") '  -0.59%  '

# Row 28
Set-TextValue $ws.Range("D28") '160.44'
Set-TextValue $ws.Range("E28") '  -1.53%  '

# Row 29
Set-TextValue $ws.Range("D29") '19.44'
Set-TextValue $ws.Range("E29") '  -1.23%  '

# Row 30
Set-TextValue $ws.Range("E30") '  +0.84%  '

# Row 31
Set-TextValue $ws.Range("D31") '1.15'
Set-TextValue $ws.Range("E31") '  -3.56%  '

# Row 32
Set-TextValue $ws.Range("E32") '  -1.18%  '

# Row 33
Set-TextValue $ws.Range("D33") '0.0619'
Set-TextValue $ws.Range("E33") '  -4.00%  '

# Row 34
Set-TextValue $ws.Range("E34") '  -0.38%  '

# Row 35
Set-TextValue $ws.Range("E35") '  -0.07%  '

# Row 36
Set-TextValue $ws.Range("B36") 'RenderToken'
Set-TextValue $ws.Range("C36") 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range("D36") '3.39'
Set-TextValue $ws.Range("E36") '  +9.72%  '

# Row 37
Set-TextValue $ws.Range("B37") 'LidoDAOToken'
Set-TextValue $ws.Range("C37") 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue $ws.Range("D37") '2.24'
Set-TextValue $ws.Range("E37") '  +1.33%  '

# Row 38
Set-TextValue $ws.Range("E38") '  -0.06%  '

# Row 39
Set-TextValue $ws.Range("D39") '5.72'
Set-TextValue $ws.Range("E39") '  -12.00%  '

# Row 40
Set-TextValue $ws.Range("D40") '0.0980'
Set-TextValue $ws.Range("E40") '  -2.09%  '

# Row 41
Set-TextValue $ws.Range("E41") '  +1.01%  '

# Row 42
Set-TextValue $ws.Range("D42") '1.17'
Set-TextValue $ws.Range("E42") '  -2.03%  '

# Row 43
Set-TextValue $ws.Range("E43") '  -0.15%  '

# Row 44
Set-TextValue $ws.Range("D44") '15.98'
Set-TextValue $ws.Range("E44") '  -3.36%  '

# Row 45
Set-TextValue $ws.Range("D45") '1.366.21'
Set-TextValue $ws.Range("E45") '  +0.21%  '

# Row 46
Set-TextValue $ws.Range("E46") '  -1.81%  '

# Row 47
Set-TextValue $ws.Range("D47") '87.85'
Set-TextValue $ws.Range("E47") '  -1.53%  '

# Row 48
Set-TextValue $ws.Range("E48") '  -2.12%  '

# Row 49
Set-TextValue $ws.Range("D49") '2.82'
Set-TextValue $ws.Range("E49") '  -0.13%  '

# Row 50
Set-TextValue $ws.Range("D50") '2.137.38'
Set-TextValue $ws.Range("E50") '  -0.23%  '

# Row 51
Set-TextValue $ws.Range("D51") '43.62'
Set-TextValue $ws.Range("E51") '  -5.74%  '
